# "updated dir and naming conventions"
#
# The R script that generates the park briefs moved out of the old
# "./scripts/" folder into "./rmd/" (matching where the other R Markdown /
# HTML assets live), and was renamed from the old "briefs" naming
# convention to the new "summaries" naming convention.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("B9").Value = "./rmd/"
$ws.Range("C9").Value = "generate_summaries_all_parks"

# Leave the selection on the cell that was being edited, matching the
# updated workbook view.
$ws.Activate()
$ws.Range("D10").Select()
